$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67.8679804978435
$ws.Range("K2").Value = 67.987520077397
$ws.Range("L2").Value = 74.0846263450239
$ws.Range("N2").Value = 64.0735542026883

$ws.Range("B3").Value = 55.2953815419129
$ws.Range("K3").Value = 56.8474042011213
$ws.Range("L3").Value = 48.4302557813727
$ws.Range("N3").Value = 48.666425997025

$ws.Range("B4").Value = 40.5493307668479
$ws.Range("K4").Value = 36.9055302021312
$ws.Range("L4").Value = 36.9722122759865
$ws.Range("N4").Value = 41.6894250824717

$ws.Range("B5").Value = 34.7872415482579
$ws.Range("K5").Value = 30.4357942185907
$ws.Range("N5").Value = 43.0029533260978

$ws.Range("B6").Value = 64.7240480131028
$ws.Range("K6").Value = 66.5933688883394
$ws.Range("L6").Value = 57.2690723237937
$ws.Range("N6").Value = 60.3471262597791

$ws.Range("B7").Value = 67.2054760884641
$ws.Range("C7").Value = 73.150166637206
$ws.Range("K7").Value = 67.292535125419
$ws.Range("L7").Value = 63.9033952437169
$ws.Range("N7").Value = 58.9468339412676

$ws.Range("B8").Value = 68.518253119246
$ws.Range("K8").Value = 70.8149579691965
$ws.Range("L8").Value = 60.4866837966034
$ws.Range("N8").Value = 63.6828257869627

$ws.Range("B9").Value = 61.895417279981
$ws.Range("K9").Value = 63.3861854984241
$ws.Range("L9").Value = 54.11746546684
$ws.Range("N9").Value = 56.0767240057917
